$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 2.25
$ws.Range("K2").Value = 2.1
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 1.9
$ws.Range("X2").Value = 15
$ws.Range("Y2").Value = 11
$ws.Range("AF2").Value = 51
$ws.Range("AL2").Value = 19
$ws.Range("AM2").Value = 29
$ws.Range("AU2").Value = 8
$ws.Range("J3").Value = 3
$ws.Range("G4").Value = 2.38
$ws.Range("I4").Value = 2.9
$ws.Range("J4").Value = 3
$ws.Range("Q4").Value = 1.92
$ws.Range("R4").Value = 1.98
$ws.Range("S4").Value = 1.36
$ws.Range("T4").Value = 3
$ws.Range("AH4").Value = 10
$ws.Range("AL4").Value = 23
$ws.Range("AT4").Value = 3
$ws.Range("S6").Value = 1.58
$ws.Range("S8").Value = 1.37
$ws.Range("H10").Value = 2.9
$ws.Range("I10").Value = 2.3
$ws.Range("S10").Value = 1.47
$ws.Range("AD10").Value = 5.5
$ws.Range("AH10").Value = 7
$ws.Range("G11").Value = 2.05
$ws.Range("I11").Value = 3.7
$ws.Range("S11").Value = 1.41
$ws.Range("T11").Value = 2.62
$ws.Range("X11").Value = 9.5
$ws.Range("Z11").Value = 19
$ws.Range("AA11").Value = 19
$ws.Range("AB11").Value = 34
$ws.Range("AH11").Value = 9.5
$ws.Range("AI11").Value = 17
$ws.Range("AO11").Value = 12
$ws.Range("AV11").Value = 51
$ws.Range("AY11").Value = 29
$ws.Range("S12").Value = 1.41
$ws.Range("T12").Value = 2.62
$ws.Range("J13").Value = 2.3
$ws.Range("L13").Value = 6
$ws.Range("O13").Value = 1.36
$ws.Range("P13").Value = 3
$ws.Range("Q13").Value = 2.2
$ws.Range("S13").Value = 1.44
$ws.Range("T13").Value = 2.63
$ws.Range("AE13").Value = 19
$ws.Range("AJ13").Value = 19
$ws.Range("AT13").Value = 2.63
$ws.Range("AW13").Value = 7
$ws.Range("R14").Value = 1.54
$ws.Range("Q15").Value = 1.54
$ws.Range("Q16").Value = 2.07
$ws.Range("R16").Value = 1.69
$ws.Range("G18").Value = 2.7
$ws.Range("I18").Value = 2.55
$ws.Range("J18").Value = 3.4
$ws.Range("L18").Value = 3.2
$ws.Range("W18").Value = 9.5
$ws.Range("AH18").Value = 9
$ws.Range("AQ18").Value = 51
$ws.Range("AW18").Value = 4.5
$ws.Range("AZ18").Value = 41
$ws.Range("M20").Value = 1.07
$ws.Range("N20").Value = 9
$ws.Range("O20").Value = 1.36
$ws.Range("P20").Value = 3
$ws.Range("G21").Value = 4.33
$ws.Range("H21").Value = 3.25
$ws.Range("I21").Value = 1.91
$ws.Range("L21").Value = 2.6
$ws.Range("N21").Value = 10
$ws.Range("Y21").Value = 15
$ws.Range("AE21").Value = 15
$ws.Range("AI21").Value = 9
$ws.Range("AJ21").Value = 8.5
$ws.Range("AL21").Value = 15
$ws.Range("AN21").Value = 6
$ws.Range("AR21").Value = 101
$ws.Range("AX21").Value = 10
$ws.Range("AZ21").Value = 34
$ws.Range("M22").Value = 1.08
$ws.Range("N22").Value = 8
$ws.Range("G24").Value = 2.3
$ws.Range("I24").Value = 3.2
$ws.Range("L24").Value = 3.6
$ws.Range("W24").Value = 8
$ws.Range("X24").Value = 11
$ws.Range("Z24").Value = 21
$ws.Range("AB24").Value = 26
$ws.Range("AJ24").Value = 12
$ws.Range("AK24").Value = 34
$ws.Range("AL24").Value = 26
$ws.Range("AP24").Value = 21
$ws.Range("AR24").Value = 51
$ws.Range("BA24").Value = 81
$ws.Range("G25").Value = 3.6
$ws.Range("H25").Value = 3.5
$ws.Range("I25").Value = 2
$ws.Range("J25").Value = 3.75
$ws.Range("L25").Value = 2.6
$ws.Range("O25").Value = 1.22
$ws.Range("P25").Value = 4
$ws.Range("U25").Value = 1.62
$ws.Range("X25").Value = 19
$ws.Range("Y25").Value = 12
$ws.Range("AA25").Value = 26
$ws.Range("AB25").Value = 29
$ws.Range("AI25").Value = 11
$ws.Range("AJ25").Value = 9
$ws.Range("AK25").Value = 19
$ws.Range("AM25").Value = 23
$ws.Range("AN25").Value = 5.5
$ws.Range("AQ25").Value = 51
$ws.Range("AW25").Value = 4.33
$ws.Range("AX25").Value = 11
$ws.Range("AZ25").Value = 34
$ws.Range("BA25").Value = 51
